# Add a new "MINUS_LOGIC_REQUIRED" config row (row 8) to the Mapping sheet,
# right after the existing UNIQUE_KEY row, pushing every following row down
# by one. This adds:
#   - A8 label "MINUS_LOGIC_REQUIRED" (same style as the other A-column
#     labels above it)
#   - B8 default value "N"
#   - a Y/N list data-validation on B8
#   - a cell comment on B8 explaining the new flag
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# Insert a new blank row at position 8 - everything currently at row 8
# and below (S.NO header, OPCO_ID.., JOIN_TABLES, LEFT rows,
# WHERE_CONDITIONS, GROUP BY, their merged cells and data validations)
# shifts down by one row automatically.
$ws.Rows.Item(8).Insert()

# Populate the new row. (Inserting the row already copied the yellow
# "label" style from row 7 above down onto row 8, matching the other
# A-column section headers such as UNIQUE_KEY / MATERIALIZATION.)
$ws.Range("A8").Value = "MINUS_LOGIC_REQUIRED"
$ws.Range("B8").Value = "N"

# Restrict B8 to Y/N via a list data validation, matching the style of
# the workbook's other list validations.
$ws.Range("B8").Validation.Add(3, 1, 1, '"Y,N"')
$ws.Range("B8").Validation.IgnoreBlank = $true
$ws.Range("B8").Validation.InCellDropdown = $true
$ws.Range("B8").Validation.ShowInput = $false
$ws.Range("B8").Validation.ShowError = $false

# Document the new flag with a cell comment.
$ws.Range("B8").AddComment("Set to Y to exclude audit columns and unique key combination in the minus logic")
